$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4155.2607
$ws.Range("I116").Value = 2079.375
$ws.Range("J116").Value = 5262.4
$ws.Range("K116").Value = 2079.375
$ws.Range("L116").Value = 5262.4
$ws.Range("M116").Value = 1362.625
$ws.Range("N116").Value = -12146.4
$ws.Range("H132").Value = 3432.8965
$ws.Range("I132").Value = 954.6667
$ws.Range("J132").Value = 9938.25
$ws.Range("K132").Value = 2864.0001
$ws.Range("L132").Value = 29814.75
$ws.Range("M132").Value = -334.0001000000002
$ws.Range("N132").Value = -34874.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2830.125
$ws.Range("I2").Value = 1113.3125
$ws.Range("J2").Value = 6263.75
$ws.Range("K2").Value = 1113.3125
$ws.Range("L2").Value = 6263.75
$ws.Range("M2").Value = -1000.3125
$ws.Range("N2").Value = -6489.75
$ws.Range("H74").Value = 1318.3721
$ws.Range("I74").Value = 899.5
$ws.Range("J74").Value = 2536.9092
$ws.Range("K74").Value = 899.5
$ws.Range("L74").Value = 2536.9092
$ws.Range("M74").Value = -25.5
$ws.Range("N74").Value = -4284.9092
$ws.Range("H77").Value = 1318.3721
$ws.Range("I77").Value = 899.5
$ws.Range("J77").Value = 2536.9092
$ws.Range("K77").Value = 4497.5
$ws.Range("L77").Value = 12684.546
$ws.Range("M77").Value = -129.5
$ws.Range("N77").Value = -21420.546
$ws.Range("H105").Value = 21900
$ws.Range("J105").Value = 21900
$ws.Range("L105").Value = 21900
$ws.Range("N105").Value = -28888
$ws.Range("H108").Value = 29900
$ws.Range("J108").Value = 29900
$ws.Range("L108").Value = 29900
$ws.Range("N108").Value = -37580
$ws.Range("H115").Value = 29950
$ws.Range("J115").Value = 29950
$ws.Range("L115").Value = 29950
$ws.Range("N115").Value = -33084
$ws.Range("H116").Value = 2830.125
$ws.Range("I116").Value = 1113.3125
$ws.Range("J116").Value = 6263.75
$ws.Range("K116").Value = 1113.3125
$ws.Range("L116").Value = 6263.75
$ws.Range("M116").Value = 1180.6875
$ws.Range("N116").Value = -10851.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2830.125
$ws.Range("I3").Value = 1113.3125
$ws.Range("J3").Value = 6263.75
$ws.Range("K3").Value = 1113.3125
$ws.Range("L3").Value = 6263.75
$ws.Range("M3").Value = -999.3125
$ws.Range("N3").Value = -6491.75
$ws.Range("H6").Value = 14758.8
$ws.Range("J6").Value = 14758.8
$ws.Range("L6").Value = 14758.8
$ws.Range("N6").Value = -14984.8
$ws.Range("H134").Value = 1771.3889
$ws.Range("I134").Value = 1737.3529
$ws.Range("J134").Value = 2350
$ws.Range("K134").Value = 5212.0587
$ws.Range("L134").Value = 7050
$ws.Range("M134").Value = -2677.0587
$ws.Range("N134").Value = -12120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1072.75
$ws.Range("I16").Value = 1027
$ws.Range("J16").Value = 1271
$ws.Range("K16").Value = 1027
$ws.Range("L16").Value = 1271
$ws.Range("M16").Value = -740
$ws.Range("N16").Value = -1845
$ws.Range("H18").Value = 37166.668
$ws.Range("J18").Value = 37166.668
$ws.Range("L18").Value = 37166.668
$ws.Range("N18").Value = -37626.668
$ws.Range("H86").Value = 40005376
$ws.Range("J86").Value = 2912.1538
$ws.Range("L86").Value = 2912.1538
$ws.Range("N86").Value = -5158.1538
$ws.Range("H89").Value = 40005376
$ws.Range("J89").Value = 2912.1538
$ws.Range("L89").Value = 14560.769
$ws.Range("N89").Value = -25792.769
$ws.Range("H113").Value = 1072.75
$ws.Range("I113").Value = 1027
$ws.Range("J113").Value = 1271
$ws.Range("K113").Value = 1027
$ws.Range("L113").Value = 1271
$ws.Range("M113").Value = 1143
$ws.Range("N113").Value = -5611
$ws.Range("H114").Value = 29830
$ws.Range("J114").Value = 29830
$ws.Range("L114").Value = 29830
$ws.Range("N114").Value = -38508
$ws.Range("H117").Value = 21765
$ws.Range("J117").Value = 21765
$ws.Range("L117").Value = 21765
$ws.Range("N117").Value = -30943

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 624.8929
$ws.Range("I11").Value = 60.4
$ws.Range("J11").Value = 938.5
$ws.Range("K11").Value = 181.2
$ws.Range("L11").Value = 2815.5
$ws.Range("M11").Value = -41.19999999999999
$ws.Range("N11").Value = -3095.5
$ws.Range("H68").Value = 2069.875
$ws.Range("J68").Value = 2312.7917
$ws.Range("L68").Value = 6938.375100000001
$ws.Range("N68").Value = -8560.375100000001
$ws.Range("H71").Value = 2069.875
$ws.Range("J71").Value = 2312.7917
$ws.Range("L71").Value = 20815.1253
$ws.Range("N71").Value = -28927.1253

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3519.1738
$ws.Range("I80").Value = 3133.4614
$ws.Range("J80").Value = 4020.6
$ws.Range("K80").Value = 3133.4614
$ws.Range("L80").Value = 4020.6
$ws.Range("M80").Value = -2135.4614
$ws.Range("N80").Value = -6016.6
$ws.Range("H83").Value = 3519.1738
$ws.Range("I83").Value = 3133.4614
$ws.Range("J83").Value = 4020.6
$ws.Range("K83").Value = 15667.307
$ws.Range("L83").Value = 20103
$ws.Range("M83").Value = -10675.307
$ws.Range("N83").Value = -30087
$ws.Range("H107").Value = 1406.75
$ws.Range("I107").Value = 1368.5714
$ws.Range("J107").Value = 1460.2
$ws.Range("K107").Value = 1368.5714
$ws.Range("L107").Value = 1460.2
$ws.Range("M107").Value = 551.4286
$ws.Range("N107").Value = -5300.2
$ws.Range("H113").Value = 1764.1111
$ws.Range("I113").Value = 1734.625
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1734.625
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 435.375
$ws.Range("N113").Value = -6340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7599.636
$ws.Range("I7").Value = 1800
$ws.Range("J7").Value = 8888.444
$ws.Range("K7").Value = 1800
$ws.Range("L7").Value = 8888.444
$ws.Range("M7").Value = -1688
$ws.Range("N7").Value = -9112.444
$ws.Range("H126").Value = 7599.636
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 8888.444
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 26665.332
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -31605.332
$ws.Range("H141").Value = 54800
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 54800
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 54800
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -65160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 45000
$ws.Range("J105").Value = 45000
$ws.Range("L105").Value = 45000
$ws.Range("N105").Value = -51988
